# Natmi following Dr Hou advice
#
# Expand the single FAPs -> Ccl21b -> Ccr10 -> sCs row into the full
# 2x2 combination of Sending cluster (FAPs, sCs) x Target cluster (sCs, ECs),
# keeping Ligand symbol (Ccl21b) and Receptor symbol (Ccr10) constant,
# and recomputing the associated NATMI statistics for each combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Ccl21b -> Ccr10 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3774076666666666
$ws.Range("H2").Value = 1.132223
$ws.Range("I2").Value = 0.6796327704557236
$ws.Range("J2").Value = 0.7608872132954309
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.5706193333333334
$ws.Range("N2").Value = 1.711858
$ws.Range("O2").Value = 0.2443278446591134
$ws.Range("P2").Value = 0.3265937887468804
$ws.Range("Q2").Value = 0.2153561111482222
$ws.Range("R2").Value = 1.938205000334
$ws.Range("S2").Value = 0.1660532099651489
$ws.Range("T2").Value = 0.2485010377992105

# Row 3: FAPs -> Ccl21b -> Ccr10 -> sCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3774076666666666
$ws.Range("H3").Value = 1.132223
$ws.Range("I3").Value = 0.6796327704557236
$ws.Range("J3").Value = 0.7608872132954309
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.7648465
$ws.Range("N3").Value = 3.529693
$ws.Range("O3").Value = 0.7556721553408867
$ws.Range("P3").Value = 0.6734062112531195
$ws.Range("Q3").Value = 0.6660665995898333
$ws.Range("R3").Value = 3.996399597539
$ws.Range("S3").Value = 0.5135795604905747
$ws.Range("T3").Value = 0.5123861754962203

# Row 4: sCs -> Ccl21b -> Ccr10 -> ECs
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.1779035
$ws.Range("H4").Value = 0.355807
$ws.Range("I4").Value = 0.3203672295442765
$ws.Range("J4").Value = 0.2391127867045691
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.5706193333333334
$ws.Range("N4").Value = 1.711858
$ws.Range("O4").Value = 0.2443278446591134
$ws.Range("P4").Value = 0.3265937887468804
$ws.Range("Q4").Value = 0.1015151765676667
$ws.Range("R4").Value = 0.6090910594060001
$ws.Range("S4").Value = 0.07827463469396451
$ws.Range("T4").Value = 0.07809275094766993

# Row 5: sCs -> Ccl21b -> Ccr10 -> sCs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.1779035
$ws.Range("H5").Value = 0.355807
$ws.Range("I5").Value = 0.3203672295442765
$ws.Range("J5").Value = 0.2391127867045691
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.7648465
$ws.Range("N5").Value = 3.529693
$ws.Range("O5").Value = 0.7556721553408867
$ws.Range("P5").Value = 0.6734062112531195
$ws.Range("Q5").Value = 0.31397236931275
$ws.Range("R5").Value = 1.255889477251
$ws.Range("S5").Value = 0.2420925948503121
$ws.Range("T5").Value = 0.1610200357568992
